# ind route: updated timeframe format
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ind")

# timeframe's accepted type changes from "int" to "string"
$ws.Range("C9").Value = "string"
# ...and now carries a note about the only supported value
$ws.Range("D9").Value = "Currently supported: H1"

# the example request URL is updated to match (options=10, timeframe=H1)
$ws.Range("C11").Value = "http://localhost:8080/ind?type=sma&options=10&pair=BTC,USD&timeframe=H1"
